$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# New header row (row 4), columns H..AA, with "GAPOK" plain header
# and two groups of colored headers (yellow = I:N, red = O:AA).
# ---------------------------------------------------------------
$headers = @{
    "H4" = "GAPOK";
    "I4" = "T. FUNGSIONAL";
    "J4" = "T. KINERJA";
    "K4" = "T. BPJS";
    "L4" = "T. STRUKTURAL";
    "M4" = "T. WALI KELAS";
    "N4" = "T. PENYESUAIAN";
    "O4" = "BPJS";
    "P4" = "Infaq TPP";
    "Q4" = "Insijam";
    "R4" = "Kalender";
    "S4" = "Koperasi/Cicilan";
    "T4" = "Lain-lain";
    "U4" = "Pinjaman Bank";
    "V4" = "Pulsa";
    "W4" = "SIMPOK";
    "X4" = "SIMWA";
    "Y4" = "Tabungan Wajib";
    "Z4" = "Verval SIMPATIKA";
    "AA4" = "Verval TPP"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# ---------------------------------------------------------------
# Data rows 5..9, columns H..AA: all zero, bordered cells.
# (Done first so the plain "border only" style lands at cellXfs
#  index 8, matching the authoring order in the real workbook.)
# ---------------------------------------------------------------
$dataRange = $ws.Range("H5:AA9")
$dataRange.Value = 0
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2

# H4 re-uses the same look as G4 (bordered, centered, no fill).
$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)

# O4:AA4 -> red fill (same border/number-format/alignment as G4).
$ws.Range("G4").Copy()
$ws.Range("O4:AA4").PasteSpecial(-4122)
$ws.Range("O4:AA4").Interior.Color = 255

# I4:N4 -> yellow fill (same border/number-format/alignment as G4).
$ws.Range("G4").Copy()
$ws.Range("I4:N4").PasteSpecial(-4122)
$ws.Range("I4:N4").Interior.Color = 65535

# ---------------------------------------------------------------
# Column widths for the new columns (approximate best-fit sizing).
# ---------------------------------------------------------------
$widths = @{
    "H" = 7.28515625; "I" = 14.7109375; "J" = 10.28515625; "K" = 7;
    "L" = 14.28515625; "M" = 13.42578125; "N" = 15.5703125; "O" = 5;
    "P" = 9.28515625; "Q" = 7.42578125; "R" = 9; "S" = 15.5703125;
    "T" = 8.5703125; "U" = 14.140625; "V" = 5.7109375; "W" = 8;
    "X" = 7.42578125; "Y" = 15.28515625; "Z" = 17; "AA" = 10.42578125
}
foreach ($col in $widths.Keys) {
    $ws.Columns($col).ColumnWidth = $widths[$col]
}

# ---------------------------------------------------------------
# View / print tweaks.
# ---------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("AA5").Select()
